# Generate Report for Handback
# Updates the "Ready for handoff" status to "Handed back: in sync with en-US"
# for the 4d1e0d08-ad0f-401d-a993-5f2f03cd48cf.md file, and records the
# handback timestamps / clears error details on the language-specific sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is 4d1e0d08-ad0f-401d-a993-5f2f03cd48cf.md
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 is 4d1e0d08-ad0f-401d-a993-5f2f03cd48cf.md
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-09-01 10:56:17"
$zhcn.Range("P3").Value = ""

# de-de sheet: row 3 is 4d1e0d08-ad0f-401d-a993-5f2f03cd48cf.md
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-09-01 10:56:24"
$dede.Range("P3").Value = ""

$zhcn.Columns.Item(16).AutoFit() | Out-Null
$dede.Columns.Item(16).AutoFit() | Out-Null
